$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibition) - update "想去人数" (want-to-go count) values in column F
$wsExhibition = $wb.Worksheets.Item("展览")
$wsExhibition.Range("F3").Value = 2136
$wsExhibition.Range("F4").Value = 878
$wsExhibition.Range("F5").Value = 1419
$wsExhibition.Range("F6").Value = 373

# Sheet "全部类型" (All Types) - same underlying events, update matching rows in column F
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F3").Value = 2136
$wsAll.Range("F6").Value = 878
$wsAll.Range("F7").Value = 1419
$wsAll.Range("F8").Value = 373
